# Daily attendance processing - 2025-11-28 03:52:43
#
# The "Recorded By" column (G) lists the users who recorded each
# attendance session as a comma-separated string (e.g. "System, someone@example.com").
# Rotate each multi-value list so the last entry moves to the front,
# leaving single-value cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ', '
    if ($parts.Count -lt 2) { continue }

    $lastPart = $parts[$parts.Count - 1]
    $rest = $parts[0..($parts.Count - 2)]
    $newVal = $lastPart + ', ' + ($rest -join ', ')

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
        $changed++
    }
}

"Recorded By column rotated on $changed rows"
